# Natmi following Dr Hou advice
#
# The LR-pair stats for Bmp4 -> Bmpr1a were recomputed with an updated
# clustering (now 3 clusters: ECs / FAPs / sCs) and now report every
# Sending-cluster x Target-cluster combination (a full 3x3 = 9-row grid
# instead of the previous 3x2 = 6 rows), so the sheet grows from
# A1:T7 to A1:T10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the various expression / specificity metrics.
$data = @(
    @("ECs",  "Bmp4", "Bmpr1a", "ECs",  2, 0.6666666666666666, 5.448315,           16.344945,          0.1199618029178375, 0.1199618029178375, 3, 1, 9.918858999999999,  29.756577,          0.160764128269069, 0.160764128269069,  54.04106827258499,  486.369614453265,   0.01928555467167201, 0.01928555467167201),
    @("ECs",  "Bmp4", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 5.448315,           16.344945,          0.1199618029178375, 0.1199618029178375, 3, 1, 29.20351433333333,  87.61054300000001,  0.47332838627826,  0.4733283862782601, 159.109945195015,   1431.989506755135,  0.05678132659013071, 0.05678132659013072),
    @("ECs",  "Bmp4", "Bmpr1a", "sCs",  2, 0.6666666666666666, 5.448315,           16.344945,          0.1199618029178375, 0.1199618029178375, 3, 1, 22.575837,           67.72751099999999,  0.3659074854526709, 0.3659074854526709, 123.000271364655,   1107.002442281895,  0.04389492165603481, 0.04389492165603482),
    @("FAPs", "Bmp4", "Bmpr1a", "ECs",  3, 1,                  31.28251333333334, 93.84754000000001,  0.6887829905701046, 0.6887829905701045, 3, 1, 9.918858999999999,  29.756577,          0.160764128269069, 0.160764128269069,  310.2868389189534,  2792.58155027058,   0.1107315970455653,  0.1107315970455652),
    @("FAPs", "Bmp4", "Bmpr1a", "FAPs", 3, 1,                  31.28251333333334, 93.84754000000001,  0.6887829905701046, 0.6887829905701045, 3, 1, 29.20351433333333,  87.61054300000001,  0.47332838627826,  0.4733283862782601, 913.5593265126913,  8222.033938614222,  0.3260205414224616,  0.3260205414224616),
    @("FAPs", "Bmp4", "Bmpr1a", "sCs",  3, 1,                  31.28251333333334, 93.84754000000001,  0.6887829905701046, 0.6887829905701045, 3, 1, 22.575837,           67.72751099999999,  0.3659074854526709, 0.3659074854526709, 706.22892196366,    6356.06029767294,   0.2520308521020777,  0.2520308521020777),
    @("sCs",  "Bmp4", "Bmpr1a", "ECs",  3, 1,                  8.686253333333333, 26.05876,            0.1912552065120579, 0.1912552065120579, 3, 1, 9.918858999999999,  29.756577,          0.160764128269069, 0.160764128269069,  86.15772205161332,  775.41949846452,    0.03074697655183176, 0.03074697655183176),
    @("sCs",  "Bmp4", "Bmpr1a", "FAPs", 3, 1,                  8.686253333333333, 26.05876,            0.1912552065120579, 0.1912552065120579, 3, 1, 29.20351433333333,  87.61054300000001,  0.47332838627826,  0.4733283862782601, 253.6691237229645,  2283.02211350668,   0.09052651826566774, 0.09052651826566775),
    @("sCs",  "Bmp4", "Bmpr1a", "sCs",  3, 1,                  8.686253333333333, 26.05876,            0.1912552065120579, 0.1912552065120579, 3, 1, 22.575837,           67.72751099999999,  0.3659074854526709, 0.3659074854526709, 196.09943939404,    1764.89495454636,   0.0699817116945584,  0.0699817116945584)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
